# Updated cryptos list (Price / Volume(1h) columns) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.999.88'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '1.626.92'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.502'
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("E8").Value = '  -1.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0618'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.38'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.80%  '
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("D12").Value = '1.853.55'
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("D14").Value = '1.620.64'
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("E15").Value = '  -2.71%  '
$ws.Range("D16").Value = '25.999.39'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '0.0₃0741'
$ws.Range("E17").Value = '  -2.32%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.79%  '
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.16%  '
$ws.Range("E21").Value = '  -2.04%  '
$ws.Range("E22").Value = '  -3.31%  '
$ws.Range("E23").Value = '  -2.02%  '
$ws.Range("E24").Value = '  +1.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("E27").Value = '  -2.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.56%  '
$ws.Range("E30").Value = '  -0.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0482'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.34%  '
$ws.Range("E32").Value = '  -3.52%  '
$ws.Range("E33").Value = '  -4.78%  '
$ws.Range("E34").Value = '  -1.32%  '
$ws.Range("E35").Value = '  -2.69%  '
$ws.Range("D36").Value = '1.123.42'
$ws.Range("E36").Value = '  -0.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.850'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.42'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.518'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.49%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0154'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.69%  '
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").Value = '1.763.49'
$ws.Range("E42").Value = '  -0.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.749'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.96%  '
$ws.Range("E45").Value = '  -1.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '54.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.41%  '
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.19%  '
$ws.Range("E51").Value = '  +0.69%  '
